# Insert a new data row at row 102 (pushing the existing rows 102..217 down
# to 103..218) and populate it with the new "Orégano" observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 6
$ws.Range("B102").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C102").Value = "Metropolitana"
$ws.Range("D102").Value = 44740
$ws.Range("E102").Value = 13
$ws.Range("F102").Value = 100112029
$ws.Range("G102").Value = "Orégano"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 34
$ws.Range("K102").Value = 14000
$ws.Range("L102").Value = 15000
$ws.Range("M102").Value = 14382
$ws.Range("N102").Value = "$/docena de atados"
$ws.Range("O102").Value = "Región Metropolitana"
$ws.Range("P102").Value = 4794
$ws.Range("Q102").Value = 3
$ws.Range("R102").Value = "Hortaliza"
